# Generate Report for Handoff
# Rows 4-7 (05b08275-..., 81931185-..., 8ecc8568-..., 94f23a9e-...) on the
# "zh-cn" and "de-de" sheets move from "low" priority to "ht" (handoff-triggered)
# priority, and their "Latest Handoff Datetime" is refreshed to reflect the
# newly generated handoff xliff files. The Overview sheet's "Latest HO Xliff
# Generate Date" column is refreshed to match the de-de handoff timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-27 08:31:26"
}

# de-de: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
for ($r = 4; $r -le 7; $r++) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-27 08:31:30"
}

# Overview: rows 4-7, Latest HO Xliff Generate Date (G) -> new timestamp
for ($r = 4; $r -le 7; $r++) {
    $overview.Range("G$r").Value = "2016-08-27 08:31:30"
}
